$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 (odds columns F:AO) ---

$row2 = @(1.58, 1.91, 4.8, 7.8, 2.88, 4.6, 1.33, 1.06, 3.15, 1.3, 1.83, 1.79, 1.32, 3, 1.84, 1.84, 1.14, 2.1, 17.5, 22, 60, 1000, 9.4, 10, 29, 110, 11, 12, 26, 1000, 19, 22, 48, 1000, 13, 1000)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 6 + $i).Value = $row2[$i]
}

$row3 = @(1.96, 2.2, 3.5, 4.2, 3.6, 4.3, 1.24, 1.05, 4.3, 1.24, 2.12, 1.7, 1.45, 2.52, 1.63, 2.24, 1.31, 1.83, 23, 21, 32, 85, 14, 11, 19, 50, 17, 12.5, 20, 55, 30, 22, 36, 90, 15, 42)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 6 + $i).Value = $row3[$i]
}

$row4 = @(2.36, 3, 3.15, 4.2, 2.66, 3.7, 1.01, 1.09, 2.66, 1.48, 1.56, 2.4, 1.2, 4.8, 1.98, 1.81, 1.33, 1.52, 1000, 1000, 1000, 1000, 9.800000000000001, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 6 + $i).Value = $row4[$i]
}

# --- Add new rows 5 and 6 ---

# Row 5: Internacional de Palmira vs Atletico Huila
$ws.Cells.Item(5, 1).Value = 'Colombian Primera B'
$ws.Range("B5").NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = '2025-11-11'
$ws.Cells.Item(5, 3).Value = '20:10:00'
$ws.Cells.Item(5, 4).Value = 'Internacional de Palmira'
$ws.Cells.Item(5, 5).Value = 'Atletico Huila'
$row5 = @(2.1, 2.98, 3, 4.9, 2.62, 980, 1.01, 1.01, 1.73, 1.01, 1.4, 2.34, 1.14, 2.32, 1.04, 1.04, 1.25, 1.51, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 6 + $i).Value = $row5[$i]
}

# Row 6: Boca Juniors de Cali vs Boyaca Patriotas
$ws.Cells.Item(6, 1).Value = 'Colombian Primera B'
$ws.Range("B6").NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = '2025-11-11'
$ws.Cells.Item(6, 3).Value = '22:20:00'
$ws.Cells.Item(6, 4).Value = 'Boca Juniors de Cali'
$ws.Cells.Item(6, 5).Value = 'Boyaca Patriotas'
$row6 = @(3.3, 4.9, 2.04, 2.68, 2.74, 3.75, 1.01, 1.1, 2.36, 1.5, 1.53, 2.3, 1.19, 5, 2.04, 1.76, 1.62, 1.25, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, 6 + $i).Value = $row6[$i]
}

